# The workbook already has the "Selection" worksheet active/selected.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Selection")

# Recollected JFreeChart (row 11) and Jodatime (row 12) manual-test results:
# previously-blank cells now record the (re-)collected counts.
$ws.Range("B11:I11").Value = 0

$ws.Range("B12:G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 0

# Move the active cell/selection to reflect where the editor ended up.
$ws.Activate()
[void]$ws.Range("G18").Select()
